# Fruta / hortaliza, semanal
# The edit reshuffles the per-row data (date + quality/volume/price/origin
# columns) among rows 2-19, while columns A-C and E-K (which are identical
# for every row in this sheet) remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row -> source row (from which the data should be copied)
$map = @{
    2  = 15
    3  = 10
    4  = 8
    5  = 19
    6  = 11
    7  = 2
    8  = 3
    9  = 13
    10 = 17
    11 = 5
    12 = 7
    13 = 4
    14 = 14
    15 = 6
    16 = 12
    17 = 18
    18 = 16
    19 = 9
}

# Columns that carry the per-row data which gets shuffled.
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)  # D, L, M, N, O, P, Q, R, S, T

# 1) Snapshot the current ("before") values for every affected column/row
#    so that writing the new values doesn't clobber a value that is still
#    needed as a source for another row.
$snapshot = @{}
foreach ($r in $map.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the new values according to the mapping.
foreach ($r in $map.Keys) {
    $src = $map[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $srcVals[$c]
    }
}
